$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7, pushing existing rows 7-9 down to 8-10.
$ws.Rows.Item(7).Insert()

# Copy the date cell style (numFmt) from the row that was just pushed down (now row 8) to the new row 7.
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat

# Fill the new row 7 with the new weekly record.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44435
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100112035
$ws.Range("G7").Value = "Bruselas (repollito)"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 140
$ws.Range("K7").Value = 21000
$ws.Range("L7").Value = 23000
$ws.Range("M7").Value = 21714
$ws.Range("N7").Value = "`$/malla 15 kilos"
$ws.Range("O7").Value = "Provincia de Quillota"
$ws.Range("P7").Value = 1448
$ws.Range("Q7").Value = 15
$ws.Range("R7").Value = "Hortaliza"
